$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$rng = $p1.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:noProof/><w:sz w:val="144"/><w:szCs w:val="144"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:t xml:space="preserve"><w:br/>[MARKER_ISOLATED 0:&lt;w:r w:rsidRPr="002C029C">&lt;w:rPr>&lt;w:noProof/>&lt;w:sz w:val="144"/>&lt;w:szCs w:val="144"/>&lt;/w:rPr>&lt;w:fldChar w:fldCharType="begin"/>&lt;/w:r>]<w:br/></w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="PMingLiU"/><w:noProof/><w:sz w:val="144"/><w:szCs w:val="144"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:t xml:space="preserve"><w:br/>[MARKER_ISOLATED 1:&lt;w:r w:rsidRPr="002C029C">&lt;w:rPr>&lt;w:rFonts w:eastAsia="PMingLiU"/>&lt;w:noProof/>&lt;w:sz w:val="144"/>&lt;w:szCs w:val="144"/>&lt;/w:rPr>&lt;w:instrText xml:space="preserve"> &lt;/w:instrText>&lt;/w:r>]<w:br/></w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="PMingLiU" w:hint="eastAsia"/><w:noProof/><w:sz w:val="144"/><w:szCs w:val="144"/></w:rPr><w:instrText>eq \o\ac(</w:instrText></w:r><w:r><w:t xml:space="preserve"><w:br/>[MARKER_ISOLATED 2:&lt;w:r w:rsidRPr="002C029C">&lt;w:rPr>&lt;w:rFonts w:eastAsia="PMingLiU" w:hint="eastAsia"/>&lt;w:noProof/>&lt;w:sz w:val="144"/>&lt;w:szCs w:val="144"/>&lt;/w:rPr>&lt;w:instrText>eq \o\ac(&lt;/w:instrText>&lt;/w:r>]<w:br/></w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="PMingLiU" w:hint="eastAsia"/><w:noProof/><w:sz w:val="144"/><w:szCs w:val="144"/></w:rPr><w:instrText>◇</w:instrText></w:r><w:r><w:t xml:space="preserve"><w:br/>[MARKER_ISOLATED 3:&lt;w:r w:rsidRPr="002C029C">&lt;w:rPr>&lt;w:rFonts w:eastAsia="PMingLiU" w:hint="eastAsia"/>&lt;w:noProof/>&lt;w:sz w:val="144"/>&lt;w:szCs w:val="144"/>&lt;/w:rPr>&lt;w:instrText>◇&lt;/w:instrText>&lt;/w:r>]<w:br/></w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="PMingLiU" w:hint="eastAsia"/><w:noProof/><w:sz w:val="144"/><w:szCs w:val="144"/></w:rPr><w:instrText>,</w:instrText></w:r><w:r><w:t xml:space="preserve"><w:br/>[MARKER_ISOLATED 4:&lt;w:r w:rsidRPr="002C029C">&lt;w:rPr>&lt;w:rFonts w:eastAsia="PMingLiU" w:hint="eastAsia"/>&lt;w:noProof/>&lt;w:sz w:val="144"/>&lt;w:szCs w:val="144"/>&lt;/w:rPr>&lt;w:instrText>,&lt;/w:instrText>&lt;/w:r>]<w:br/></w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="PMingLiU" w:hint="eastAsia"/><w:noProof/><w:position w:val="4"/><w:sz w:val="144"/><w:szCs w:val="144"/></w:rPr><w:instrText>M</w:instrText></w:r><w:r><w:t xml:space="preserve"><w:br/>[MARKER_ISOLATED 5:&lt;w:r w:rsidRPr="002C029C">&lt;w:rPr>&lt;w:rFonts w:eastAsia="PMingLiU" w:hint="eastAsia"/>&lt;w:noProof/>&lt;w:position w:val="4"/>&lt;w:sz w:val="144"/>&lt;w:szCs w:val="144"/>&lt;/w:rPr>&lt;w:instrText>M&lt;/w:instrText>&lt;/w:r>]<w:br/></w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="PMingLiU" w:hint="eastAsia"/><w:noProof/><w:sz w:val="144"/><w:szCs w:val="144"/></w:rPr><w:instrText>)</w:instrText></w:r><w:r><w:t xml:space="preserve"><w:br/>[MARKER_ISOLATED 6:&lt;w:r w:rsidRPr="002C029C">&lt;w:rPr>&lt;w:rFonts w:eastAsia="PMingLiU" w:hint="eastAsia"/>&lt;w:noProof/>&lt;w:sz w:val="144"/>&lt;w:szCs w:val="144"/>&lt;/w:rPr>&lt;w:instrText>)&lt;/w:instrText>&lt;/w:r>]<w:br/></w:t></w:r><w:r><w:rPr><w:noProof/><w:sz w:val="144"/><w:szCs w:val="144"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r><w:t xml:space="preserve"><w:br/>[MARKER_ISOLATED 7:&lt;w:r w:rsidRPr="002C029C">&lt;w:rPr>&lt;w:noProof/>&lt;w:sz w:val="144"/>&lt;w:szCs w:val="144"/>&lt;/w:rPr>&lt;w:fldChar w:fldCharType="end"/>&lt;/w:r>]<w:br/></w:t></w:r></w:p>
'@

$rng.InsertXML($xml)
"done"
